$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.589.34"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "'1.923.54"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'247.15"
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4732"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("D8").Value = "'0.2917"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "'0.06851"
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("D10").Value = "'106.00"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").Value = "'18.51"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "'1.933.77"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "'5.350"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").Value = "'0.6744"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "'289.17"
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("D17").Value = "'30.615.18"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'0.000007655"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'5.573"
$ws.Range("E21").Value = "  +5.30%  "
$ws.Range("D22").Value = "'2.181.46"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'6.493"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").Value = "'9.559"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").Value = "'166.94"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "'20.83"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "'2.136"
$ws.Range("E28").Value = "  +3.76%  "
$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").Value = "'1.407"
$ws.Range("E30").Value = "  +3.31%  "
$ws.Range("D31").Value = "'4.211"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("D32").Value = "'4.081"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").Value = "'0.05055"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").Value = "'0.7355"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "'0.02057"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("D37").Value = "'2.745"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'0.9997"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "'2.057"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'111.67"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").Value = "'0.4470"
$ws.Range("E42").Value = "  +6.40%  "
$ws.Range("D43").Value = "'0.8739"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "'5.921"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "'68.15"
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("D47").Value = "'7.334"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'49.12"
$ws.Range("E48").Value = "  +13.01%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.416"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "'0.1266"
$ws.Range("E50").Value = "  +4.04%  "
$ws.Range("D51").Value = "'35.35"
$ws.Range("E51").Value = "  +1.26%  "
